$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Films")
$h = $ws.Range("F5").Hyperlinks.Item(1)
Write-Host ("count before: " + $ws.Hyperlinks.Count)
$h.Delete()
Write-Host ("count after direct delete: " + $ws.Hyperlinks.Count)
